$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the "Save" header in H1, copying format from the adjacent "sum" header (G1)
# so it picks up the same bold/centered/bordered style used by the rest of row 1.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Fill in the "Save" column values for rows 2-13
$saveValues = @(0, 0, 0, 0, 1, 1, 0, 1, 0, 0, 0, 0)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
